# Apply the target edits to the "MASTERFILE EQUIPMENT" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix equipment description for the first equipment (row 8)
$ws.Range("D8").Value = "Air Receiver"

# 2) Normalize pressure unit formatting ("4 Bar.G" -> "4 Bar G",
#    "1 Bar.G" -> "1 Bar G") for every data row (8 through 35)
for ($row = 8; $row -le 35; $row++) {
    $ws.Cells.Item($row, 13).Value = "4 Bar G"   # column M: DESIGN PRESSURE
    $ws.Cells.Item($row, 15).Value = "1 Bar G"   # column O: OPERATING PRESSURE
}

# 3) For rows 26-35, the material "TYPE" (column H) could not be verified,
#    so mark it "Not Found", and clear out the now-unknown "GRADE" (column J)
for ($row = 26; $row -le 35; $row++) {
    $ws.Cells.Item($row, 8).Value = "Not Found"
    $ws.Range("J$row").ClearContents()
}
